$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Insert a new row above row 1, shifting existing data down by one row
$ws.Rows.Item(1).Insert()

# Put the new header value into the freshly inserted A1
$ws.Range("A1").Value = "Columna A"

# Move the selection to C3, matching the saved view state
$ws.Range("C3").Select()
